# Applies the Oct 20 2023 cryptos-list price/volume refresh to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into column D while forcing it to stay text.
# Excel auto-converts plain numeric-looking strings like "27.01" into
# real numbers, which would lose trailing zeros / reformat them. We
# use a leading apostrophe to force text entry, then reset the cell
# style back to "Normal" so no stray number-format style is left behind.
function Set-TextCell($cell, $text) {
    $cell.Value = "`'" + $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '29.636.52'
$ws.Range("E2").Value = '  +3.54%  '
$ws.Range("D3").Value = '1.609.18'
$ws.Range("E3").Value = '  +2.83%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("E5").Value = '  +1.15%  '
Set-TextCell $ws.Range("D6") '0.520'
$ws.Range("E6").Value = '  +2.23%  '
$ws.Range("E7").Value = '  +0.16%  '
Set-TextCell $ws.Range("D8") '27.01'
$ws.Range("E8").Value = '  +8.45%  '
Set-TextCell $ws.Range("D9") '43.61'
$ws.Range("E9").Value = '  -1.28%  '
$ws.Range("E10").Value = '  +2.49%  '
Set-TextCell $ws.Range("D11") '0.0601'
$ws.Range("E11").Value = '  +2.50%  '
$ws.Range("E12").Value = '  +1.42%  '
$ws.Range("D13").Value = '1.839.15'
$ws.Range("E13").Value = '  +2.86%  '
$ws.Range("D14").Value = '1.655.42'
$ws.Range("E14").Value = '  +5.87%  '
$ws.Range("D15").Value = '29.660.15'
$ws.Range("E15").Value = '  +3.53%  '
$ws.Range("E16").Value = '  +3.88%  '
$ws.Range("E17").Value = '  +2.58%  '
Set-TextCell $ws.Range("D18") '63.54'
$ws.Range("E18").Value = '  +3.41%  '
Set-TextCell $ws.Range("D19") '240.30'
$ws.Range("E19").Value = '  +5.60%  '
Set-TextCell $ws.Range("D20") '7.60'
$ws.Range("E20").Value = '  +3.82%  '
$ws.Range("D21").Value = '0.0₃0694'
$ws.Range("E21").Value = '  +1.90%  '
$ws.Range("E22").Value = '  +0.07%  '
Set-TextCell $ws.Range("D23") '3.99'
$ws.Range("E23").Value = '  +1.58%  '
Set-TextCell $ws.Range("D24") '9.23'
$ws.Range("E24").Value = '  +1.86%  '
Set-TextCell $ws.Range("D25") '2.10'
$ws.Range("E25").Value = '  +0.90%  '
Set-TextCell $ws.Range("D26") '154.85'
$ws.Range("E26").Value = '  +2.01%  '
Set-TextCell $ws.Range("D27") '15.31'
$ws.Range("E27").Value = '  +3.51%  '
Set-TextCell $ws.Range("D28") '0.108'
$ws.Range("E28").Value = '  +2.57%  '
$ws.Range("E29").Value = '  +2.80%  '
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("E31").Value = '  +3.53%  '
$ws.Range("E32").Value = '  +0.97%  '
Set-TextCell $ws.Range("D33") '3.22'
$ws.Range("E33").Value = '  +1.19%  '
$ws.Range("E34").Value = '  +4.33%  '
$ws.Range("D35").Value = '1.419.55'
$ws.Range("E35").Value = '  +1.00%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("E37").Value = '  +4.86%  '
$ws.Range("E38").Value = '  +5.34%  '
$ws.Range("E39").Value = '  +0.34%  '
$ws.Range("E40").Value = '  +2.34%  '
$ws.Range("E41").Value = '  +4.65%  '
Set-TextCell $ws.Range("D42") '1.98'
$ws.Range("E42").Value = '  +2.21%  '
Set-TextCell $ws.Range("D43") '54.46'
$ws.Range("E43").Value = '  +27.82%  '
$ws.Range("E44").Value = '  +6.57%  '
Set-TextCell $ws.Range("D45") '0.800'
$ws.Range("E45").Value = '  +4.14%  '
$ws.Range("E46").Value = '  +0.09%  '
Set-TextCell $ws.Range("D47") '65.96'
$ws.Range("E47").Value = '  +3.22%  '
Set-TextCell $ws.Range("D48") '5.30'
$ws.Range("E48").Value = '  +1.43%  '
$ws.Range("D49").Value = '1.750.86'
$ws.Range("E49").Value = '  +3.16%  '
Set-TextCell $ws.Range("D50") '0.874'
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("E51").Value = '  +2.29%  '
